$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'65.623.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -1.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'3.446.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -3.85%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'595.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -1.62%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'136.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  -7.85%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'3.445.33"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -3.87%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'  +0.58%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'7.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  -4.88%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Formula = "'  -9.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Formula = "'  -8.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'4.032.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  -3.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -11.22%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'3.470.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -3.14%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'26.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  -10.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'65.567.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  -1.27%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Formula = "'  -2.24%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'9.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -10.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Formula = "'  -8.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'13.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -7.45%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'394.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -6.53%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -10.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'73.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  -6.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Formula = "'  -0.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'3.593.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -3.77%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'0.0000105"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  -11.64%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  -0.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Formula = "'  -10.11%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Formula = "'  -8.92%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'8.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  -12.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'3.455.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  -3.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Formula = "'  +0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Formula = "'  -7.37%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Formula = "'22.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  -8.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'173.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  -1.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Formula = "'  -13.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Formula = "'  -10.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  -8.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'4.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  -13.62%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.0780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  -8.46%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.818"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -6.93%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'43.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  -5.37%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +0.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'4.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -14.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'1.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  -11.81%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'23.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -3.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'1.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -2.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'  -8.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  -16.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'2.210.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  -7.53%  "
$ws.Range("E51").Style = "Normal"
